# Versionshinweise_Zeugnisse_BK.docx update
# Add a new "17.01.2022" changelog entry right after the introductory
# paragraph "Versionshinweise zu den Zeugnissen:" at the very top of the
# document, consisting of a date line followed by a bulleted note about
# the new B7 Kinderpfleger certificate wording.

$d = $word.ActiveDocument

# The very first paragraph in the document is the intro line.
$introPara = $d.Paragraphs.Item(1)

# --- Insert the date paragraph "17.01.2022" --------------------------------
$dateRange = $introPara.Range
$dateRange.Collapse(0)
$dateRange.InsertParagraphAfter()

$dateRange = $d.Paragraphs.Item(2).Range
$dateRange.Text = "17.01.2022"

# --- Insert the new bulleted changelog entry --------------------------------
$dateRange = $d.Paragraphs.Item(2).Range
$dateRange.Collapse(0)
$dateRange.InsertParagraphAfter()

$itemPara = $d.Paragraphs.Item(3)
$itemRange = $itemPara.Range
$itemRange.Text = "Änderung des Satzes zur Zusatzqualifikation nach QHB auf dem ASZ B7 für die Kinderpfleger. Umgesetzt auf einem neuen Formular."

# Give it the same "Listenabsatz" bullet-list look used by all the other
# changelog bullet points in this document.
$itemPara.Style = "Listenabsatz"

$gallery = $word.ListGalleries.Item(1)
$template = $gallery.ListTemplates.Item(1)
$itemPara.Range.ListFormat.ApplyListTemplateWithLevel($template, $false, 1, $false, 0)

Write-Host "Inserted 17.01.2022 changelog entry."
